$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 18: recalibrated raw temperature readings ("rectal" row).
$ws.Range("B18").Value = -4.57
$ws.Range("C18").Value = 0.26
$ws.Range("D18").Value = 5.13
$ws.Range("E18").Value = 9.77
$ws.Range("F18").Value = 10.07
$ws.Range("G18").Value = 20.02
$ws.Range("H18").Value = 29.99
$ws.Range("I18").Value = 39.97
$ws.Range("J18").Value = 50

# D18 switches from the orange highlight (s=4, used by E18/F18) to the red
# highlight used by B18/C18 (s=3) -- copy the font formatting over.
$ws.Range("D18").Font.Color = $ws.Range("C18").Font.Color

# Restore the view: scrolled down a bit further and a different active cell.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("K37").Select()
